# Weekly update: a new Alcachofa (Hortaliza, Macroferia Regional de Talca)
# price record is inserted at the top of the data block (row 67), pushing
# the existing records (old rows 67-90) down by one row (new rows 68-91).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 67; this shifts old rows 67:90 down to 68:91
# and extends the used range from A1:R90 to A1:R91 automatically.
$ws.Rows.Item(67).Insert()

# Populate the newly-inserted row 67 with this week's new record.
$ws.Cells.Item(67, 1).Value  = 5
$ws.Cells.Item(67, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(67, 3).Value  = "Maule"
$ws.Cells.Item(67, 4).Value  = 44782
$ws.Cells.Item(67, 5).Value  = 7
$ws.Cells.Item(67, 6).Value  = 100112013
$ws.Cells.Item(67, 7).Value  = "Alcachofa"
$ws.Cells.Item(67, 8).Value  = "Madrigal"
$ws.Cells.Item(67, 9).Value  = "Primera"
$ws.Cells.Item(67, 10).Value = 300
$ws.Cells.Item(67, 11).Value = 14000
$ws.Cells.Item(67, 12).Value = 14000
$ws.Cells.Item(67, 13).Value = 14000
$ws.Cells.Item(67, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(67, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(67, 16).Value = 350
$ws.Cells.Item(67, 17).Value = 40
$ws.Cells.Item(67, 18).Value = "Hortaliza"
